$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number;
# force Text format first, assign, then restore default style so no
# extraneous number-format styling remains on the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4764'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4047'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08476'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.060'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.39'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.599'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.188'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001072'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06626'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.012'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.867'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.897'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.167'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '124.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9825'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09641'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.456'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.702'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.623'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.153'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02334'
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6211'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1917'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.351'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5949'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.058'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.413'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06814'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.27'
$ws.Range("D51").Style = "Normal"

# Plain text assignments (safe from numeric auto-detection already).
$ws.Range("D2").Value = '28.529.43'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").Value = '1.964.98'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -4.16%  '
$ws.Range("E8").Value = '  -3.58%  '
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("E10").Value = '  -6.45%  '
$ws.Range("E11").Value = '  -3.12%  '
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("D13").Value = '1.985.54'
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("E14").Value = '  -3.10%  '
$ws.Range("E15").Value = '  -3.47%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("E19").Value = '  -0.66%  '
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").Value = '28.573.45'
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = '2.227.77'
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("E29").Value = '  -4.78%  '
$ws.Range("E30").Value = '  -3.88%  '
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("E32").Value = '  -5.18%  '
$ws.Range("E33").Value = '  -1.87%  '
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -2.99%  '
$ws.Range("E37").Value = '  +1.99%  '
$ws.Range("E38").Value = '  -3.50%  '
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("E41").Value = '  -3.26%  '
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("E44").Value = '  -3.81%  '
$ws.Range("E45").Value = '  +5.48%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("E47").Value = '  -3.76%  '
$ws.Range("E48").Value = '  -4.98%  '
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("E51").Value = '  -0.97%  '
